$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "test #"
$ws.Range("B1").Value = "Species "
$ws.Range("C1").Value = "Mooring"
$ws.Range("D1").Value = "detTotal"
$ws.Range("E1").Value = "numTP"
$ws.Range("F1").Value = "numFP"
$ws.Range("G1").Value = "numFN"
$ws.Range("H1").Value = "TPR"
$ws.Range("I1").Value = "FPR"
$ws.Range("J1").Value = "TPdivFP"
$ws.Range("K1").Value = "AUC"
$ws.Range("M1").Value = "Total counts"
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "GS"
$ws.Range("C2").Value = "rf AW15_AU_BS3_files_705-749"
$ws.Range("D2").Value = 958
$ws.Range("E2").Value = 633
$ws.Range("F2").Value = 325
$ws.Range("G2").Value = 64
$ws.Range("H2").Value = 0.90817790530846498
$ws.Range("I2").Value = 0.339248434237996
$ws.Range("J2").Value = 1.9476923076923101
$ws.Range("K2").Value = 0.94405257587029101
$ws.Range("M2").Value = 2169
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = "GS"
$ws.Range("C3").Value = "rf BS12_AU_02a_files_1-46"
$ws.Range("D3").Value = 508
$ws.Range("E3").Value = 390
$ws.Range("F3").Value = 118
$ws.Range("G3").Value = 11
$ws.Range("H3").Value = 0.972568578553616
$ws.Range("I3").Value = 0.232283464566929
$ws.Range("J3").Value = 3.3050847457627102
$ws.Range("K3").Value = 0.94405257587029101
$ws.Range("M3").Value = 758
$ws.Range("A4").Value = 0
$ws.Range("B4").Value = "GS"
$ws.Range("C4").Value = "rf AW14_AU_BS3_files_1-71"
$ws.Range("D4").Value = 1453
$ws.Range("E4").Value = 944
$ws.Range("F4").Value = 509
$ws.Range("G4").Value = 75
$ws.Range("H4").Value = 0.92639842983317
$ws.Range("I4").Value = 0.350309704060564
$ws.Range("J4").Value = 1.8546168958742599
$ws.Range("K4").Value = 0.94405257587029101
$ws.Range("M4").Value = 3187
$ws.Range("A5").Value = 0
$ws.Range("B5").Value = "GS"
$ws.Range("C5").Value = "rf BS13_AU_04_files_137-224"
$ws.Range("D5").Value = 1320
$ws.Range("E5").Value = 654
$ws.Range("F5").Value = 666
$ws.Range("G5").Value = 35
$ws.Range("H5").Value = 0.94920174165457205
$ws.Range("I5").Value = 0.50454545454545496
$ws.Range("J5").Value = 0.98198198198198205
$ws.Range("K5").Value = 0.94405257587029101
$ws.Range("M5").Value = 3660
$ws.Range("A6").Value = 0
$ws.Range("B6").Value = "GS"
$ws.Range("C6").Value = "rf AW12_AU_BS3_files_1-250"
$ws.Range("D6").Value = 2632
$ws.Range("E6").Value = 1619
$ws.Range("F6").Value = 1013
$ws.Range("G6").Value = 188
$ws.Range("H6").Value = 0.895960154952961
$ws.Range("I6").Value = 0.38487841945288798
$ws.Range("J6").Value = 1.59822309970385
$ws.Range("K6").Value = 0.94405257587029101
$ws.Range("M6").Value = 5560
$ws.Range("A7").Value = 0
$ws.Range("B7").Value = "GS"
$ws.Range("C7").Value = "rf AW12_AU_BS3_files_1464-1507"
$ws.Range("D7").Value = 728
$ws.Range("E7").Value = 552
$ws.Range("F7").Value = 176
$ws.Range("G7").Value = 67
$ws.Range("H7").Value = 0.89176090468497604
$ws.Range("I7").Value = 0.24175824175824201
$ws.Range("J7").Value = 3.1363636363636398
$ws.Range("K7").Value = 0.94405257587029101
$ws.Range("M7").Value = 1126
$ws.Range("A8").Value = 0
$ws.Range("B8").Value = "GS"
$ws.Range("C8").Value = "rf all"
$ws.Range("D8").Value = 7599
$ws.Range("E8").Value = 4792
$ws.Range("F8").Value = 2807
$ws.Range("G8").Value = 440
$ws.Range("H8").Value = 0.91590214067278297
$ws.Range("I8").Value = 0.36939070930385598
$ws.Range("J8").Value = 1.7071606697541899
$ws.Range("K8").Value = 0.94405257587029101
$ws.Range("M8").Value = 16460
$ws.Range("A9").Value = "average"
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "GS"
$ws.Range("C10").Value = "rf AW14_AU_BS3_files_309-369"
$ws.Range("D10").Value = "."
$ws.Range("E10").Value = "."
$ws.Range("F10").Value = "."
$ws.Range("G10").Value = "."
$ws.Range("H10").Value = 0.91590214067278297
$ws.Range("I10").Value = 0.36939070930385598
$ws.Range("J10").Value = 1.7071606697541899
$ws.Range("K10").Value = 0.94405257587029101
$ws.Range("M10").Value = 3344
$ws.Range("A11").Value = "actual"
$ws.Range("A12").Value = 1
$ws.Range("B12").Value = "GS"
$ws.Range("C12").Value = "rf AW14_AU_BS3_files_309-369"
$ws.Range("D12").Value = 909
$ws.Range("E12").Value = 626
$ws.Range("F12").Value = 283
$ws.Range("G12").Value = 92
$ws.Range("H12").Value = 0.871866295264624
$ws.Range("I12").Value = 0.31133113311331101
$ws.Range("J12").Value = 2.2120141342756199
$ws.Range("K12").Value = 0.88861900000000005
$ws.Range("M12").Value = 3344
$ws.Range("A13").Value = "Combined"
$ws.Range("A14").Value = 1
$ws.Range("B14").Value = "GS"
$ws.Range("C14").Value = "rf AW15_AU_BS3_files_705-749"
$ws.Range("D14").Value = 1000
$ws.Range("E14").Value = 639
$ws.Range("F14").Value = 361
$ws.Range("G14").Value = 58
$ws.Range("H14").Value = 0.91678622668579601
$ws.Range("I14").Value = 0.36099999999999999
$ws.Range("J14").Value = 1.7700831024930701
$ws.Range("K14").Value = 0.94777387726579698
$ws.Range("M14").Value = 2169
$ws.Range("A15").Value = 1
$ws.Range("B15").Value = "GS"
$ws.Range("C15").Value = "rf BS12_AU_02a_files_1-46"
$ws.Range("D15").Value = 516
$ws.Range("E15").Value = 393
$ws.Range("F15").Value = 123
$ws.Range("G15").Value = 8
$ws.Range("H15").Value = 0.98004987531172105
$ws.Range("I15").Value = 0.23837209302325599
$ws.Range("J15").Value = 3.1951219512195101
$ws.Range("K15").Value = 0.94777387726579698
$ws.Range("M15").Value = 758
$ws.Range("A16").Value = 1
$ws.Range("B16").Value = "GS"
$ws.Range("C16").Value = "rf AW14_AU_BS3_files_1-71"
$ws.Range("D16").Value = 1489
$ws.Range("E16").Value = 952
$ws.Range("F16").Value = 537
$ws.Range("G16").Value = 67
$ws.Range("H16").Value = 0.93424926398429797
$ws.Range("I16").Value = 0.36064472800537301
$ws.Range("J16").Value = 1.77281191806331
$ws.Range("K16").Value = 0.94777387726579698
$ws.Range("M16").Value = 3187
$ws.Range("A17").Value = 1
$ws.Range("B17").Value = "GS"
$ws.Range("C17").Value = "rf BS13_AU_04_files_137-224"
$ws.Range("D17").Value = 1416
$ws.Range("E17").Value = 670
$ws.Range("F17").Value = 746
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 0.97242380261248196
$ws.Range("I17").Value = 0.52683615819208995
$ws.Range("J17").Value = 0.898123324396783
$ws.Range("K17").Value = 0.94777387726579698
$ws.Range("M17").Value = 3660
$ws.Range("A18").Value = 1
$ws.Range("B18").Value = "GS"
$ws.Range("C18").Value = "rf AW12_AU_BS3_files_1-250"
$ws.Range("D18").Value = 2736
$ws.Range("E18").Value = 1631
$ws.Range("F18").Value = 1105
$ws.Range("G18").Value = 176
$ws.Range("H18").Value = 0.90260099612617595
$ws.Range("I18").Value = 0.403874269005848
$ws.Range("J18").Value = 1.4760180995475101
$ws.Range("K18").Value = 0.94777387726579698
$ws.Range("M18").Value = 5560
$ws.Range("A19").Value = 1
$ws.Range("B19").Value = "GS"
$ws.Range("C19").Value = "rf AW12_AU_BS3_files_1464-1507"
$ws.Range("D19").Value = 739
$ws.Range("E19").Value = 554
$ws.Range("F19").Value = 185
$ws.Range("G19").Value = 65
$ws.Range("H19").Value = 0.89499192245557302
$ws.Range("I19").Value = 0.25033829499323401
$ws.Range("J19").Value = 2.9945945945945902
$ws.Range("K19").Value = 0.94777387726579698
$ws.Range("M19").Value = 1126
$ws.Range("A20").Value = 1
$ws.Range("B20").Value = "GS"
$ws.Range("C20").Value = "rf AW14_AU_BS3_files_309-369"
$ws.Range("D20").Value = 859
$ws.Range("E20").Value = 614
$ws.Range("F20").Value = 245
$ws.Range("G20").Value = 104
$ws.Range("H20").Value = 0.85515320334261802
$ws.Range("I20").Value = 0.28521536670547099
$ws.Range("J20").Value = 2.5061224489795899
$ws.Range("K20").Value = 0.94777387726579698
$ws.Range("M20").Value = 3344
$ws.Range("A21").Value = 1
$ws.Range("B21").Value = "GS"
$ws.Range("C21").Value = "rf all"
$ws.Range("D21").Value = 8755
$ws.Range("E21").Value = 5453
$ws.Range("F21").Value = 3302
$ws.Range("G21").Value = 497
$ws.Range("H21").Value = 0.91647058823529404
$ws.Range("I21").Value = 0.377155910908053
$ws.Range("J21").Value = 1.6514233797698401
$ws.Range("K21").Value = 0.94777387726579698
$ws.Range("M21").Formula = "=SUM(M14:M20)"
$ws.Range("A22").Value = "average"
$ws.Range("A23").Value = 2
$ws.Range("B23").Value = "GS"
$ws.Range("C23").Value = "GS AW15_AU_BS2_files_33-103"
$ws.Range("D23").Value = "."
$ws.Range("E23").Value = "."
$ws.Range("F23").Value = "."
$ws.Range("G23").Value = "."
$ws.Range("H23").Value = 0.91647058823529404
$ws.Range("I23").Value = 0.377155910908053
$ws.Range("J23").Value = 1.6514233797698401
$ws.Range("K23").Value = 0.94777387726579698
$ws.Range("M23").Value = 3282
$ws.Range("A24").Value = "actual"
$ws.Range("A25").Value = 2
$ws.Range("B25").Value = "GS"
$ws.Range("C25").Value = "GS AW15_AU_BS2_files_33-103"
$ws.Range("D25").Value = 956
$ws.Range("E25").Value = 416
$ws.Range("F25").Value = 540
$ws.Range("G25").Value = 64
$ws.Range("H25").Value = 0.86666666666666703
$ws.Range("I25").Value = 0.56485355648535596
$ws.Range("J25").Value = 0.77037037037037004
$ws.Range("K25").Value = 0.84250133547008499
$ws.Range("M25").Value = 3282
$ws.Range("A26").Value = "Combined"
$ws.Range("A27").Value = 2
$ws.Range("B27").Value = "GS"
$ws.Range("C27").Value = "rf AW15_AU_BS3_files_705-749"
$ws.Range("D27").Value = 994
$ws.Range("E27").Value = 641
$ws.Range("F27").Value = 353
$ws.Range("G27").Value = 56
$ws.Range("H27").Value = 0.91965566714490699
$ws.Range("I27").Value = 0.35513078470824899
$ws.Range("J27").Value = 1.8158640226628899
$ws.Range("K27").Value = 0.948787406201965
$ws.Range("M27").Value = 2169
$ws.Range("A28").Value = 2
$ws.Range("B28").Value = "GS"
$ws.Range("C28").Value = "rf BS12_AU_02a_files_1-46"
$ws.Range("D28").Value = 529
$ws.Range("E28").Value = 394
$ws.Range("F28").Value = 135
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = 0.98254364089775603
$ws.Range("I28").Value = 0.25519848771266501
$ws.Range("J28").Value = 2.9185185185185198
$ws.Range("K28").Value = 0.948787406201965
$ws.Range("M28").Value = 758
$ws.Range("A29").Value = 2
$ws.Range("B29").Value = "GS"
$ws.Range("C29").Value = "rf AW14_AU_BS3_files_1-71"
$ws.Range("D29").Value = 1507
$ws.Range("E29").Value = 952
$ws.Range("F29").Value = 555
$ws.Range("G29").Value = 67
$ws.Range("H29").Value = 0.93424926398429797
$ws.Range("I29").Value = 0.36828135368281401
$ws.Range("J29").Value = 1.71531531531532
$ws.Range("K29").Value = 0.948787406201965
$ws.Range("M29").Value = 3187
$ws.Range("A30").Value = 2
$ws.Range("B30").Value = "GS"
$ws.Range("C30").Value = "rf BS13_AU_04_files_137-224"
$ws.Range("D30").Value = 1498
$ws.Range("E30").Value = 678
$ws.Range("F30").Value = 820
$ws.Range("G30").Value = 11
$ws.Range("H30").Value = 0.98403483309143702
$ws.Range("I30").Value = 0.54739652870494004
$ws.Range("J30").Value = 0.826829268292683
$ws.Range("K30").Value = 0.948787406201965
$ws.Range("M30").Value = 3660
$ws.Range("A31").Value = 2
$ws.Range("B31").Value = "GS"
$ws.Range("C31").Value = "rf AW12_AU_BS3_files_1-250"
$ws.Range("D31").Value = 2884
$ws.Range("E31").Value = 1648
$ws.Range("F31").Value = 1236
$ws.Range("G31").Value = 159
$ws.Range("H31").Value = 0.91200885445489799
$ws.Range("I31").Value = 0.42857142857142899
$ws.Range("J31").Value = 1.3333333333333299
$ws.Range("K31").Value = 0.948787406201965
$ws.Range("M31").Value = 5560
$ws.Range("A32").Value = 2
$ws.Range("B32").Value = "GS"
$ws.Range("C32").Value = "rf AW12_AU_BS3_files_1464-1507"
$ws.Range("D32").Value = 755
$ws.Range("E32").Value = 557
$ws.Range("F32").Value = 198
$ws.Range("G32").Value = 62
$ws.Range("H32").Value = 0.89983844911147004
$ws.Range("I32").Value = 0.26225165562913899
$ws.Range("J32").Value = 2.81313131313131
$ws.Range("K32").Value = 0.948787406201965
$ws.Range("M32").Value = 1126
$ws.Range("A33").Value = 2
$ws.Range("B33").Value = "GS"
$ws.Range("C33").Value = "rf AW14_AU_BS3_files_309-369"
$ws.Range("D33").Value = 865
$ws.Range("E33").Value = 619
$ws.Range("F33").Value = 246
$ws.Range("G33").Value = 99
$ws.Range("H33").Value = 0.86211699164345401
$ws.Range("I33").Value = 0.284393063583815
$ws.Range("J33").Value = 2.5162601626016299
$ws.Range("K33").Value = 0.948787406201965
$ws.Range("M33").Value = 3344
$ws.Range("A34").Value = 2
$ws.Range("B34").Value = "GS"
$ws.Range("C34").Value = "rf AW15_AU_BS2_files_33-103"
$ws.Range("D34").Value = 750
$ws.Range("E34").Value = 398
$ws.Range("F34").Value = 352
$ws.Range("G34").Value = 76
$ws.Range("H34").Value = 0.83966244725738404
$ws.Range("I34").Value = 0.46933333333333299
$ws.Range("J34").Value = 1.1306818181818199
$ws.Range("K34").Value = 0.948787406201965
$ws.Range("M34").Value = 3281
$ws.Range("A35").Value = 2
$ws.Range("B35").Value = "GS"
$ws.Range("C35").Value = "rf all"
$ws.Range("D35").Value = 9782
$ws.Range("E35").Value = 5887
$ws.Range("F35").Value = 3895
$ws.Range("G35").Value = 537
$ws.Range("H35").Value = 0.91640722291407195
$ws.Range("I35").Value = 0.39818033122060897
$ws.Range("J35").Value = 1.5114249037227201
$ws.Range("K35").Value = 0.948787406201965
$ws.Range("M35").Formula = "=SUM(M27:M34)"

$ws.Range("R32").Select()
